# The deck ships two DrawingML themes:
#   ppt/theme/theme1.xml -> attached to the slide master (color scheme "Integral")
#   ppt/theme/theme2.xml -> attached to the notes master  (color scheme "Office")
#
# The commit swaps the two color schemes: the slide master (theme1.xml) now
# carries the stock "Office" palette, while the notes master (theme2.xml)
# picks up the old "Integral" palette. Font scheme and format scheme are
# identical between the two themes, so only the 12 theme colors change.
#
# Helper: turn a "RRGGBB" hex string into the packed BGR integer that
# PowerPoint's RGB()/ColorFormat.RGB expects (0x00BBGGRR).
function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Target palette for the slide master's theme (theme1.xml): the stock
# "Office" color scheme, in the standard dk1/lt1/dk2/lt2/accent1-6/hlink/
# folHlink order used by ThemeColorScheme.Item(1..12).
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = HexToComRgb($officeColors[$i - 1])
}
